# paper edits, responding to Shaul
#
# 1) The auto-updating "date" placeholder field (type datetimeFigureOut)
#    on the slide master and every slide layout moves from 5/18/2018 to
#    6/3/2018 (this is PowerPoint re-stamping the auto date field on save).
# 2) On the slide itself, two callout labels get their numbers tweaked:
#       "Primary mirror 15 K"      -> "Primary mirror 17 K"   (3 runs)
#       "Low pass filter 100 mK"   -> "Low pass filter 1 K"   (2 runs)

$p = $ppt.ActivePresentation

# ppPlaceholderDate
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $shp.TextFrame.TextRange.Text = "6/3/2018"
        }
    }
}

# --- Slide master date placeholder ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- Every slide layout's date placeholder ---
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $lyt = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $lyt.Shapes
}

# --- Content edits on slide 1 ---
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $shp = $grp.GroupItems.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    if (-not $shp.TextFrame.HasText) { continue }
    $txt = $shp.TextFrame.TextRange.Text

    if ($txt -eq "Primary mirror 15 K") {
        $tr = $shp.TextFrame.TextRange
        $tr.Text = "Primary mirror "
        $tr.InsertAfter("17 ")
        $tr.InsertAfter("K")
    }
    elseif ($txt -eq "Low pass filter 100 mK") {
        $tr = $shp.TextFrame.TextRange
        $tr.Text = "Low pass filter "
        $tr.InsertAfter("1 K")
    }
}
